# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml  -> currently "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml  -> currently "Integral"     (used by the Slide Master / all slides)
#
# The authored change swaps the two themes: the design that slides/Slide Master
# use becomes the default "Office Theme" palette, while the (otherwise unused)
# "Integral" palette is what's left associated with the Notes Master side.
#
# The only part of that swap reachable from the slide-facing object model is the
# 12-colour theme colour scheme that every slide (via its layout/master) shares,
# so we repoint those 12 colours from the "Integral" palette to the stock
# "Office Theme" palette.

function Get-BgrInt([string]$hex) {
    # PowerPoint's RGB-valued COM properties pack colours as 0x00BBGGRR
    # (the same layout VBA's RGB() function produces), not 0x00RRGGBB.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in the fixed clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = Get-BgrInt $officeTheme[$i - 1]
}
